{"js": "// 1) Update the letter date from September 19, 2025 to September 21, 2025.\nconst dateResults = context.document.body.search(\"September 19, 2025\", { matchCase: true, matchWholeWord: false });\ndateResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"September 21, 2025\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Split the mailing address paragraph \"3380 Eichers Pl, Santa Clara CA 95051\"\n//    (the standalone address line in the letterhead, not the one inside the\n//    account-summary table) into two paragraphs:\n//      \"3380 Eichers Pl\"\n//      \"Santa Clara, CA 95051\"\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text === \"3380 Eichers Pl, Santa Clara CA 95051\") {\n    para.load(\"parentTableOrNullObject\");\n    await context.sync();\n    if (para.parentTableOrNullObject.isNullObject) {\n      para.insertText(\"3380 Eichers Pl\", \"Replace\");\n      para.insertParagraph(\"Santa Clara, CA 95051\", \"After\");\n      await context.sync();\n      break;\n    }\n  }\n}\n\n// 3) Remove the now-redundant empty \"No Spacing\" paragraph that immediately\n//    follows the \"...Board of Directors\" signature line.\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items/text,items/style\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text === \"Townhomes at Nuevo Homeowners Association Board of Directors\") {\n    const next = paragraphs2.items[i + 1];\n    if (next) {\n      next.load(\"text,style\");\n      await context.sync();\n      if (next.text === \"\" && next.style === \"No Spacing\") {\n        next.delete();\n        await context.sync();\n      }\n    }\n    break;\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the letter date from September 19, 2025 to September 21, 2025.\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"September 19, 2025\") {\n        $p.Range.Text = \"September 21, 2025\"\n        break\n    }\n}\n\n# 2) Split the mailing address paragraph \"3380 Eichers Pl, Santa Clara CA 95051\"\n#    (the standalone address line in the letterhead, not the copy of it that\n#    lives inside the account-summary table) into two paragraphs:\n#      \"3380 Eichers Pl\"\n#      \"Santa Clara, CA 95051\"\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"3380 Eichers Pl, Santa Clara CA 95051\" -and $p.Range.Information(12) -eq $false) {\n        $r = $p.Range\n        $r.Text = \"3380 Eichers Pl\"\n        $r.InsertParagraphAfter()\n        $newPara = $d.Paragraphs.Item($i + 1)\n        $newPara.Range.Text = \"Santa Clara, CA 95051\"\n        break\n    }\n}\n\n# 3) Remove the now-redundant empty \"No Spacing\" paragraph that immediately\n#    follows the \"...Board of Directors\" signature line.\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"Townhomes at Nuevo Homeowners Association Board of Directors\") {\n        $next = $d.Paragraphs.Item($i + 1)\n        $nextText = $next.Range.Text.TrimEnd([char]13, [char]7)\n        if ($nextText -eq \"\" -and $next.Style.NameLocal -eq \"No Spacing\") {\n            $next.Range.Delete()\n        }\n        break\n    }\n}\n"}
